# Fruta / hortaliza, semanal
# Inserts the new weekly data points (week of 2021-09-21) for
# Femacal de La Calera - Piña at row 210-211, pushing all subsequent
# historical rows down by two (old row 210 -> 212, ... old row 312 -> 314).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 210.
$ws.Rows.Item(210).EntireRow.Insert()
$ws.Rows.Item(210).EntireRow.Insert()

# New weekly records (same shape/order as the existing columns A..T).
$rowsData = @(
  @(3, "Femacal de La Calera", "Coquimbo", 44460, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Primera", 162, 20000, 20000, 20000, "`$/caja 12 unidades", "Ecuador", 1667, 12),
  @(3, "Femacal de La Calera", "Coquimbo", 44460, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Segunda", 108, 20000, 20000, 20000, "`$/caja 14 unidades", "Ecuador", 1429, 14)
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
  $r = 210 + $i
  $data = $rowsData[$i]
  for ($j = 0; $j -lt $data.Length; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $data[$j]
  }
}
